$d = $word.ActiveDocument

# The document contains one table. Its last column holds per-student
# totals such as "+0" that must be normalized to "0" (the leading
# "+" sign is dropped). Other cells such as "3+3" must be left untouched,
# so walk each cell explicitly and only touch ones whose full text is
# exactly "+0".
$t = $d.Tables.Item(1)
$lastCol = $t.Columns.Count

for ($r = 1; $r -le $t.Rows.Count; $r++) {
    try {
        $c = $t.Cell($r, $lastCol)
    } catch {
        continue
    }
    $txt = $c.Range.Text
    $txt = $txt.TrimEnd([char]7, [char]13)
    if ($txt -eq "+0") {
        $c.Range.Text = "0"
    }
}
